$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.948.56'
$ws.Range('E2').Value = '  -2.41%  '
$ws.Range('D3').Value = '1.893.63'
$ws.Range('E3').Value = '  -5.56%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.27%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '324.09'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.45%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.53%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4596'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -2.45%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.3820'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -3.64%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '45.59'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -3.13%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.07736'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -2.92%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.9701'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -3.75%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '22.12'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D13').Value = '1.903.40'
$ws.Range('E13').Value = '  -6.15%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '5.691'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '6.962'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -4.34%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.07047'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('E17').Value = '  -0.46%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '83.61'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -6.34%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.000009545'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -4.63%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '16.72'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -4.02%  '
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').Value = '28.883.04'
$ws.Range('E22').Value = '  -2.88%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.320'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -4.30%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '10.92'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -3.37%  '
$ws.Range('D25').Value = '2.115.96'
$ws.Range('E25').Value = '  -7.04%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.068'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -2.93%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '156.14'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -1.73%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '19.06'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -3.61%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '5.607'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -6.31%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '117.58'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -2.66%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.818'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -8.01%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '0.09265'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -2.40%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.8543'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -4.96%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '5.087'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -4.21%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.240'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -7.95%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '3.013'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -5.81%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.05693'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -2.86%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '1.148'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -2.70%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '1.001'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.46%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.02042'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -4.58%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '7.432'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -6.39%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.5510'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -4.80%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.1752'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -4.21%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.000002870'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -16.84%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '9.252'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -6.56%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.698'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +1.35%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.5185'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -4.13%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '11.29'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -7.47%  '
$ws.Range('E49').Value = '  -4.76%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.06808'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -2.75%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '111.67'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -2.71%  '
